# Update countries & provincias Spain
# - Refresh the "last updated" timestamp (04:04 -> 05:04)
# - Refresh COVID case numbers for a handful of countries, which also
#   changes their ranking order (so the country names in column A for
#   rows 80/81, 153/154 and 192/193 swap places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp row
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 05:04"

# Row 53 - Australia (no re-rank, only refreshed figures)
$ws.Range("B53").Value = 6939
$ws.Range("C53").Value = 10
$ws.Range("D53").Value = 6141
$ws.Range("E53").Value = 701

# Rows 80/81 - Honduras overtakes Islandia
$ws.Range("A80").Value = "Honduras"
$ws.Range("B80").Value = 1830
$ws.Range("C80").Value = 59
$ws.Range("D80").Value = 195
$ws.Range("E80").Value = 1527
$ws.Range("F80").Value = 10
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 108

$ws.Range("A81").Value = "Islandia"
$ws.Range("B81").Value = 1801
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 1773
$ws.Range("E81").Value = 18
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 10

# Rows 153/154 - Guayana Francesa overtakes Republica de Africa Central
$ws.Range("A153").Value = "Guayana Francesa"
$ws.Range("B153").Value = 144
$ws.Range("C153").Value = 3
$ws.Range("D153").Value = 122
$ws.Range("E153").Value = 21
$ws.Range("H153").Value = 1

$ws.Range("A154").Value = "Republica de Africa Central"
$ws.Range("B154").Value = 143
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 10
$ws.Range("E154").Value = 133
$ws.Range("H154").Value = 0

# Row 161 - Nepal (no re-rank, only refreshed figures)
$ws.Range("B161").Value = 110
$ws.Range("C161").Value = 1
$ws.Range("E161").Value = 79

# Rows 192/193 - Belice overtakes Nueva Caledonia
$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0
